{"js": "// Convert legacy `<w:fldSimple w:instr=\"...\">` field codes into the\n// equivalent \"complex field\" run sequence that Word itself writes once a\n// document containing a `fldSimple` has been opened/edited and saved again:\n//\n//   <w:r><w:fldChar w:fldCharType=\"begin\"/></w:r>\n//   <w:r><w:instrText>INSTR</w:instrText></w:r>\n//   <w:r><w:fldChar w:fldCharType=\"separate\"/></w:r>\n//   <w:r><w:fldChar w:fldCharType=\"end\"/></w:r>\n//\n// This mirrors the M2Doc fix: the simple-field form is technically valid\n// OOXML but some downstream tooling only understands the expanded form, so\n// the content is otherwise lost after a round trip.\n\nfunction escapeXmlText(s) {\n  return s.replace(/&/g, \"&amp;\").replace(/</g, \"&lt;\").replace(/>/g, \"&gt;\");\n}\n\nfunction decodeXmlAttr(s) {\n  return s\n    .replace(/&quot;/g, '\"')\n    .replace(/&apos;/g, \"'\")\n    .replace(/&lt;/g, \"<\")\n    .replace(/&gt;/g, \">\")\n    .replace(/&amp;/g, \"&\");\n}\n\nfunction buildFldCharRuns(instr) {\n  return (\n    '<w:r><w:fldChar w:fldCharType=\"begin\"/></w:r>' +\n    \"<w:r><w:instrText>\" + escapeXmlText(instr) + \"</w:instrText></w:r>\" +\n    '<w:r><w:fldChar w:fldCharType=\"separate\"/></w:r>' +\n    '<w:r><w:fldChar w:fldCharType=\"end\"/></w:r>'\n  );\n}\n\n// `Paragraph.getOoxml()` synthesizes `w14:paraId`/`w14:textId` attributes\n// that are not present in the authored document - drop them so we do not\n// introduce attributes that were never in the original file.\nfunction stripSyntheticParaIds(pTagOpen) {\n  return pTagOpen\n    .replace(/\\s+w14:paraId=\"[^\"]*\"/, \"\")\n    .replace(/\\s+w14:textId=\"[^\"]*\"/, \"\");\n}\n\n// Replace every `<w:fldSimple w:instr=\"...\">...</w:fldSimple>` (or the\n// self-closed form) found in a paragraph's XML with the expanded\n// begin/instrText/separate/end run sequence, leaving everything else\n// (paragraph attributes, bookmarks, other runs, etc.) untouched.\nfunction expandFldSimple(paragraphXml) {\n  const fldSimpleRe = /<w:fldSimple\\s+w:instr=\"([^\"]*)\"\\s*(?:\\/>|>[\\s\\S]*?<\\/w:fldSimple>)/g;\n  return paragraphXml.replace(fldSimpleRe, (whole, instr) => {\n    return buildFldCharRuns(decodeXmlAttr(instr));\n  });\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Find every paragraph that currently hosts at least one field (the only\n// fields in this document are the legacy `fldSimple` ones).\nconst paraFieldRanges = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const range = paragraphs.items[i].getRange();\n  range.load(\"fields\");\n  paraFieldRanges.push(range);\n}\nawait context.sync();\n\nconst targetIndexes = [];\nfor (let i = 0; i < paraFieldRanges.length; i++) {\n  if (paraFieldRanges[i].fields.items.length > 0) {\n    targetIndexes.push(i);\n  }\n}\n\nif (targetIndexes.length > 0) {\n  // Grab the exact authored XML (attributes, bookmarks, etc. included) for\n  // every paragraph we are about to touch.\n  const ooxmlRequests = {};\n  for (const idx of targetIndexes) {\n    ooxmlRequests[idx] = paragraphs.items[idx].getOoxml();\n  }\n  await context.sync();\n\n  for (const idx of targetIndexes) {\n    const fullXml = ooxmlRequests[idx].value;\n    const bodyMatch = fullXml.match(/<w:body>([\\s\\S]*)<\\/w:body>/);\n    if (!bodyMatch) {\n      continue;\n    }\n    // The paragraph we asked for is always the first element serialized\n    // into the synthetic single-paragraph package body.\n    const pMatch = bodyMatch[1].match(/^<w:p(?:\\s[^>]*)?>[\\s\\S]*?<\\/w:p>|^<w:p(?:\\s[^>]*)?\\/>/);\n    if (!pMatch) {\n      continue;\n    }\n    let paragraphXml = stripSyntheticParaIds(pMatch[0]);\n    if (!/fldSimple/.test(paragraphXml)) {\n      continue;\n    }\n    const newParagraphXml = expandFldSimple(paragraphXml);\n\n    const packageXml =\n      '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n      '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n      '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      \"<pkg:xmlData>\" +\n      '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n      \"<w:body>\" + newParagraphXml + \"</w:body>\" +\n      \"</w:document>\" +\n      \"</pkg:xmlData>\" +\n      \"</pkg:part>\" +\n      \"</pkg:package>\";\n\n    paragraphs.items[idx].insertOoxml(packageXml, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Convert legacy `<w:fldSimple w:instr=\"...\">` field codes into the\n# equivalent \"complex field\" run sequence that Word itself writes once a\n# document containing a `fldSimple` has been opened/edited and saved again:\n#\n#   <w:r><w:fldChar w:fldCharType=\"begin\"/></w:r>\n#   <w:r><w:instrText>INSTR</w:instrText></w:r>\n#   <w:r><w:fldChar w:fldCharType=\"separate\"/></w:r>\n#   <w:r><w:fldChar w:fldCharType=\"end\"/></w:r>\n#\n# This mirrors the M2Doc fix: the simple-field form is technically valid\n# OOXML but some downstream tooling only understands the expanded form, so\n# the content is otherwise lost after a round trip.\n\n$d = $word.ActiveDocument\n\n# `<w:fldSimple w:instr=\"...\">...</w:fldSimple>` either self-closed or with\n# an explicit close tag (the engine always round-trips it with an explicit\n# close tag, but stay tolerant of both forms).\n$fldSimpleRegex = [regex]'(?s)<w:fldSimple\\s+w:instr=\"([^\"]*)\"\\s*(?:/>|>.*?</w:fldSimple>)'\n\n# `Paragraph.Range.WordOpenXML` synthesizes `w14:paraId`/`w14:textId`\n# attributes that are not present in the authored document - drop them so we\n# do not introduce attributes that were never in the original file.\n$paraIdRegex = [regex]'\\s+w14:paraId=\"[^\"]*\"'\n$textIdRegex = [regex]'\\s+w14:textId=\"[^\"]*\"'\n\n# Find every paragraph that currently hosts at least one field (the only\n# fields in this document are the legacy `fldSimple` ones).\n$targetIndexes = New-Object System.Collections.ArrayList\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    if ($para.Range.Fields.Count -gt 0) {\n        [void]$targetIndexes.Add($i)\n    }\n}\n\nforeach ($idx in $targetIndexes) {\n    $para = $d.Paragraphs.Item($idx)\n    $range = $para.Range\n\n    # Recover the exact authored XML (attributes, bookmarks, etc. included)\n    # for this paragraph only.\n    $fullXml = $range.WordOpenXML\n\n    $bodyStart = $fullXml.IndexOf(\"<w:body>\")\n    $bodyEnd = $fullXml.IndexOf(\"</w:body>\")\n    if ($bodyStart -lt 0 -or $bodyEnd -lt 0) {\n        continue\n    }\n    $bodyInner = $fullXml.Substring($bodyStart + 8, $bodyEnd - ($bodyStart + 8))\n\n    # The paragraph we asked for is always the first element serialized into\n    # the synthetic single-paragraph package body.\n    $pRegex = [regex]'(?s)^<w:p(?:\\s[^>]*)?>.*?</w:p>|^<w:p(?:\\s[^>]*)?/>'\n    $pMatch = $pRegex.Match($bodyInner)\n    if (-not $pMatch.Success) {\n        continue\n    }\n    $paragraphXml = $pMatch.Value\n    $paragraphXml = $paraIdRegex.Replace($paragraphXml, \"\")\n    $paragraphXml = $textIdRegex.Replace($paragraphXml, \"\")\n\n    if ($paragraphXml -notmatch \"fldSimple\") {\n        continue\n    }\n\n    # Expand every `fldSimple` found in this paragraph into the begin /\n    # instrText / separate / end run sequence, leaving everything else\n    # (paragraph attributes, bookmarks, other runs, etc.) untouched.\n    $newParagraphXml = $paragraphXml\n    while ($true) {\n        $m = $fldSimpleRegex.Match($newParagraphXml)\n        if (-not $m.Success) {\n            break\n        }\n        $instr = $m.Groups[1].Value\n        $replacement = ('<w:r><w:fldChar w:fldCharType=\"begin\"/></w:r>' +\n            '<w:r><w:instrText>' + $instr + '</w:instrText></w:r>' +\n            '<w:r><w:fldChar w:fldCharType=\"separate\"/></w:r>' +\n            '<w:r><w:fldChar w:fldCharType=\"end\"/></w:r>')\n        $newParagraphXml = $newParagraphXml.Substring(0, $m.Index) + $replacement + $newParagraphXml.Substring($m.Index + $m.Length)\n    }\n\n    $packageXml = ('<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $newParagraphXml + '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>')\n\n    [void]$range.InsertXML($packageXml)\n}\n"}
